# Natmi following Dr Hou advice
# Update the Edn1-Ednra ligand-receptor pair table to include the "ECs"
# sending/target cluster alongside the existing "FAPs" and "sCs" clusters,
# and refresh the computed statistics for the full 3x3 cluster matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Edn1","Ednra","ECs",2,0.6666666666666666,7.803861666666666,23.411585,0.7185921299200277,0.7185921299200276,3,1,2.614605666666666,7.843817,0.07238665237615237,0.07238665237615237,20.40402093554944,183.636188419945,0.05201647870875997,0.05201647870875996),
    @("ECs","Edn1","Ednra","FAPs",2,0.6666666666666666,7.803861666666666,23.411585,0.7185921299200277,0.7185921299200276,3,1,14.70158366666666,44.10475099999999,0.4070206225838464,0.4070206225838464,114.7291252155928,1032.562126940335,0.2924818161039019,0.2924818161039019),
    @("ECs","Edn1","Ednra","sCs",2,0.6666666666666666,7.803861666666666,23.411585,0.7185921299200277,0.7185921299200276,3,1,18.80380766666667,56.411423,0.5205927250400012,0.5205927250400012,146.7423138372728,1320.680824535455,0.3740938351073658,0.3740938351073658),
    @("FAPs","Edn1","Ednra","ECs",3,1,2.696109,8.088327,0.2482620517329206,0.2482620517329206,3,1,2.614605666666666,7.843817,0.07238665237615237,0.07238665237615237,7.049261869350999,63.443356824159,0.01797085883698128,0.01797085883698128),
    @("FAPs","Edn1","Ednra","FAPs",3,1,2.696109,8.088327,0.2482620517329206,0.2482620517329206,3,1,14.70158366666666,44.10475099999999,0.4070206225838464,0.4070206225838464,39.63707203795299,356.733648341577,0.1010477748602764,0.1010477748602764),
    @("FAPs","Edn1","Ednra","sCs",3,1,2.696109,8.088327,0.2482620517329206,0.2482620517329206,3,1,18.80380766666667,56.411423,0.5205927250400012,0.5205927250400012,50.697115084369,456.274035759321,0.1292434180356629,0.1292434180356629),
    @("sCs","Edn1","Ednra","ECs",3,1,0.3599613333333334,1.079884,0.03314581834705166,0.03314581834705165,3,1,2.614605666666666,7.843817,0.07238665237615237,0.07238665237615237,0.9411569419142222,8.470412477228001,0.002399314830411122,0.002399314830411121),
    @("sCs","Edn1","Ednra","FAPs",3,1,0.3599613333333334,1.079884,0.03314581834705166,0.03314581834705165,3,1,14.70158366666666,44.10475099999999,0.4070206225838464,0.4070206225838464,5.292001658764889,47.62801492888399,0.01349103161966804,0.01349103161966804),
    @("sCs","Edn1","Ednra","sCs",3,1,0.3599613333333334,1.079884,0.03314581834705166,0.03314581834705165,3,1,18.80380766666667,56.411423,0.5205927250400012,0.5205927250400012,6.76864367943689,60.917793114932,0.01725547189697249,0.01725547189697249)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
